$wb = $excel.ActiveWorkbook

# The workbook has two sheets that duplicate the same event listing:
#   "展览"     (sheet1) - rows 2-16
#   "全部类型" (sheet4) - rows 2-19
# Update the "想去人数" (F column) counters on both sheets.

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 93
$ws1.Range("F3").Value = 4046
$ws1.Range("F4").Value = 2367
$ws1.Range("F8").Value = 28
$ws1.Range("F13").Value = 1511
$ws1.Range("F15").Value = 2876

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 93
$ws4.Range("F3").Value = 4046
$ws4.Range("F4").Value = 2367
$ws4.Range("F8").Value = 28
$ws4.Range("F16").Value = 1511
$ws4.Range("F18").Value = 2876
